$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New issue #11 appended as the next row in the issues table.
# The "Issue ID" column holds text values elsewhere in the sheet
# (e.g. A2:A6 are the text "6".."10"), so force this one to text
# too via the leading apostrophe, matching how Excel stores typed
# text that looks like a number.
$ws.Range("A7").Value = "'11"
$ws.Range("B7").Value = "Resource newrelic_synthetics_monitor: Cannot unset validation_string"
$ws.Range("C7").Value = "open"
$ws.Range("D7").Value = "2025-03-24T09:02:25Z"
$ws.Range("E7").Value = "bug"
